$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 410
$ws.Range("J6").Value = 682.6667
$ws.Range("N6").Value = -2272.0001
$ws.Range("L6").Value = 2048.0001
$ws.Range("N8").Value = -6164
$ws.Range("J8").Value = 1962
$ws.Range("I8").Value = 22.5
$ws.Range("H8").Value = 551.4545000000001
$ws.Range("L8").Value = 5886
$ws.Range("M8").Value = 71.5
$ws.Range("K8").Value = 67.5
$ws.Range("N33").Value = -547.4
$ws.Range("I33").Value = 141
$ws.Range("K33").Value = 141
$ws.Range("J33").Value = 89.40000000000001
$ws.Range("H33").Value = 124.875
$ws.Range("M33").Value = 88
$ws.Range("L33").Value = 89.40000000000001
$ws.Range("K98").Value = 1090
$ws.Range("H98").Value = 1064.1428
$ws.Range("I98").Value = 1090
$ws.Range("M98").Value = 408
$ws.Range("J116").Value = 5516
$ws.Range("I116").Value = 9166.5
$ws.Range("K116").Value = 9166.5
$ws.Range("H116").Value = 7341.25
$ws.Range("N116").Value = -12400
$ws.Range("M116").Value = -5724.5
$ws.Range("L116").Value = 5516
$ws.Range("M122").Value = -820
$ws.Range("H122").Value = 1064.1428
$ws.Range("K122").Value = 3270
$ws.Range("I122").Value = 1090
$ws.Range("K135").Value = 44380.287
$ws.Range("M135").Value = -41845.287
$ws.Range("H135").Value = 5854.1
$ws.Range("J135").Value = 8007.6665
$ws.Range("N135").Value = -77138.9985
$ws.Range("L135").Value = 72068.9985
$ws.Range("I135").Value = 4931.143
$ws.Range("M137").Value = -2773.5
$ws.Range("H137").Value = 1774.5
$ws.Range("I137").Value = 1774.5
$ws.Range("K137").Value = 5323.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M4").Value = -209
$ws.Range("H4").Value = 325
$ws.Range("I4").Value = 325
$ws.Range("K4").Value = 325
$ws.Range("I110").Value = 5286485.5
$ws.Range("H110").Value = 3701099.8
$ws.Range("J110").Value = 1866
$ws.Range("K110").Value = 5286485.5
$ws.Range("N110").Value = -5956
$ws.Range("L110").Value = 1866
$ws.Range("M110").Value = -5284440.5
$ws.Range("M122").Value = -1534
$ws.Range("N122").Value = -16898.5
$ws.Range("H122").Value = 2218.5
$ws.Range("L122").Value = 11998.5
$ws.Range("K122").Value = 3984
$ws.Range("J122").Value = 3999.5
$ws.Range("I122").Value = 1328

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M134").Value = -2693.571599999999
$ws.Range("I134").Value = 1742.8572
$ws.Range("H134").Value = 1742.8572
$ws.Range("K134").Value = 5228.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M22").Value = 250
$ws.Range("H22").Value = 100
$ws.Range("K22").Value = 100
$ws.Range("I22").Value = 100
$ws.Range("N31").Value = -4737.6665
$ws.Range("M31").Value = -1696
$ws.Range("I31").Value = 1991
$ws.Range("K31").Value = 1991
$ws.Range("J31").Value = 4147.6665
$ws.Range("L31").Value = 4147.6665
$ws.Range("H31").Value = 2709.889
$ws.Range("L34").Value = 4147.6665
$ws.Range("N34").Value = -4551.6665
$ws.Range("H34").Value = 2709.889
$ws.Range("K34").Value = 1991
$ws.Range("I34").Value = 1991
$ws.Range("M34").Value = -1789
$ws.Range("J34").Value = 4147.6665
$ws.Range("L60").Value = 0
$ws.Range("H60").Value = 7000
$ws.Range("N60").ClearContents()
$ws.Range("J60").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J7").Value = 1739.6666
$ws.Range("I7").Value = 186.33333
$ws.Range("K7").Value = 558.99999
$ws.Range("H7").Value = 963
$ws.Range("L7").Value = 5218.9998
$ws.Range("M7").Value = -446.99999
$ws.Range("N7").Value = -5442.9998
$ws.Range("H12").Value = 121
$ws.Range("L12").Value = 173.727276
$ws.Range("J12").Value = 57.909092
$ws.Range("N12").Value = -519.7272760000001
$ws.Range("M15").Value = -2116
$ws.Range("J15").Value = 401
$ws.Range("H15").Value = 518
$ws.Range("L15").Value = 1203
$ws.Range("I15").Value = 752
$ws.Range("N15").Value = -1483
$ws.Range("K15").Value = 2256
$ws.Range("K23").Value = 180000270
$ws.Range("H23").Value = 30000334
$ws.Range("I23").Value = 60000090
$ws.Range("J23").Value = 578
$ws.Range("M23").Value = -180000035
$ws.Range("L23").Value = 1734
$ws.Range("N23").Value = -2204
$ws.Range("H98").Value = 895.6667
$ws.Range("N98").Value = -5757.875
$ws.Range("J98").Value = 920.625
$ws.Range("L98").Value = 2761.875
$ws.Range("H131").Value = 502030.06
$ws.Range("N131").Value = -2517836.82
$ws.Range("J131").Value = 835918.9399999999
$ws.Range("L131").Value = 2507756.82
$ws.Range("L137").Value = 11301.6
$ws.Range("N137").Value = -21501.6
$ws.Range("J137").Value = 3767.2
$ws.Range("M137").Value = 2625.75
$ws.Range("H137").Value = 2459.4443
$ws.Range("I137").Value = 824.75
$ws.Range("K137").Value = 2474.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 15566
$ws.Range("N98").Value = -21556
$ws.Range("J98").Value = 15566
$ws.Range("L98").Value = 15566
$ws.Range("M122").Value = -11094.1432
$ws.Range("N122").Value = -7900
$ws.Range("H122").Value = 3733.6667
$ws.Range("L122").Value = 3000
$ws.Range("K122").Value = 13544.1432
$ws.Range("J122").Value = 1000
$ws.Range("I122").Value = 4514.7144
$ws.Range("K132").Value = 25135.251
$ws.Range("N132").Value = -24879.5
$ws.Range("M132").Value = -22605.251
$ws.Range("I132").Value = 8378.416999999999
$ws.Range("L132").Value = 19819.5
$ws.Range("J132").Value = 6606.5
$ws.Range("H132").Value = 8125.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J7").Value = 8936.75
$ws.Range("I7").Value = 30001
$ws.Range("K7").Value = 30001
$ws.Range("H7").Value = 11277.223
$ws.Range("L7").Value = 8936.75
$ws.Range("M7").Value = -29889
$ws.Range("N7").Value = -9160.75
$ws.Range("I40").Value = 1865.6666
$ws.Range("K40").Value = 1865.6666
$ws.Range("H40").Value = 3816.1667
$ws.Range("L40").Value = 5766.6665
$ws.Range("J40").Value = 5766.6665
$ws.Range("N40").Value = -6038.6665
$ws.Range("M40").Value = -1729.6666
$ws.Range("M46").Value = -619.75
$ws.Range("I46").Value = 807.75
$ws.Range("L46").Value = 1974.8334
$ws.Range("H46").Value = 1508
$ws.Range("J46").Value = 1974.8334
$ws.Range("K46").Value = 807.75
$ws.Range("N46").Value = -2350.8334
$ws.Range("H126").Value = 11277.223
$ws.Range("I126").Value = 30001
$ws.Range("N126").Value = -31750.25
$ws.Range("J126").Value = 8936.75
$ws.Range("L126").Value = 26810.25
$ws.Range("M126").Value = -87533
$ws.Range("K126").Value = 90003
$ws.Range("H138").Value = 90000
$ws.Range("J138").Value = 90000
$ws.Range("L138").Value = 90000
$ws.Range("N138").Value = -100280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2404.75
$ws.Range("I126").Value = 1522.5714
$ws.Range("N126").Value = -15859.4
$ws.Range("J126").Value = 3639.8
$ws.Range("L126").Value = 10919.4
$ws.Range("M126").Value = -2097.7142
$ws.Range("K126").Value = 4567.7142
$ws.Range("K132").Value = 10274.25
$ws.Range("M132").Value = -7744.25
$ws.Range("I132").Value = 3424.75
$ws.Range("H132").Value = 3424.75
